$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B19 to be a true numeric value (was stored as text/inline string)
$ws.Range("B19").Value = 3

# Add new row 20 with the additional annotation record
$ws.Range("A20").Value = "Sunsi Wu"

# B20 must stay a text value ("5"), not be auto-converted to a number.
# Type it with a leading apostrophe into a scratch cell, then copy/paste
# just the value so the destination keeps the plain (unstyled) cell format.
$ws.Range("ZZ1").Value = "'5"
$ws.Range("ZZ1").Copy()
$ws.Range("B20").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("C20").Value = "No clear novelty"
$ws.Range("D20").Value = "CRT"
$ws.Range("E20").Value = "MET"
$ws.Range("F20").Value = "4efacd8b-a5d8-471d-9660-f5eb687b96fc"
$ws.Range("G20").Value = "Byni8NLHf_annotated.xlsx"
$ws.Range("H20").Value = "No clear novelty"
